$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.228.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.783.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.551"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0657"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.039.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.774.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("E15").Value = "  -3.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.224.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.93%  "

$ws.Range("E24").Value = "  -3.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("E32").Value = "  -2.03%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.440.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.23%  "

$ws.Range("E36").Value = "  -2.31%  "

$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.66%  "

$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("E42").Value = "  -3.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.64%  "

$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.940.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.36%  "
